$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1) Rows 15 and 16 had their match data (columns F:V) swapped - the
#    "Sukhothai vs Trat FC" match now sits on row 15 and the
#    "Buriram vs Lamphun Warrior" match now sits on row 16 (columns A:E,
#    which hold the row index/country/tournament/season/date, stay put).
# -----------------------------------------------------------------------
$row15New = @("Sukhothai", 0, "Trat FC", 0, 1.78, "18/08/2023 23:26", 1.83, "20/08/2023 12:51", 3.66, "18/08/2023 23:26", 3.87, "20/08/2023 12:51", 4.14, "18/08/2023 23:26", 4.1, "20/08/2023 12:51", "https://www.betexplorer.com/football/thailand/thai-league-1/sukhothai-trat-fc/xGYJjFEC/")
$row16New = @("Buriram", 3, "Lamphun Warrior", 0, 1.25, "15/08/2023 16:42", 1.21, "20/08/2023 12:50", 6.02, "15/08/2023 16:42", 6.56, "20/08/2023 12:50", 10.71, "15/08/2023 16:42", 12.97, "20/08/2023 12:50", "https://www.betexplorer.com/football/thailand/thai-league-1/buriram-united-f-c-lamphun-warrior/nLdbew0s/")

for ($i = 0; $i -lt $row15New.Count; $i++) {
    $col = 6 + $i
    $ws.Cells.Item(15, $col).Value = $row15New[$i]
    $ws.Cells.Item(16, $col).Value = $row16New[$i]
}

# -----------------------------------------------------------------------
# 2) Append 11 new match rows (110-120) at the bottom of the sheet.
#    Each inner array holds columns A..V in order.
# -----------------------------------------------------------------------
$newRows = @(
    @(109, "thailand", "thai-league-1", "2023-2024", 45283.5, "Chonburi", 2, "Muang Thong Utd", 1, 1.76, "16/12/2023 12:43", 1.98, "23/12/2023 11:56", 3.93, "16/12/2023 12:43", 4.08, "23/12/2023 11:56", 3.97, "16/12/2023 12:43", 3.37, "23/12/2023 11:56", "https://www.betexplorer.com/football/thailand/thai-league-1/chonburi-muang-thong-utd/zu6jhs87/"),
    @(110, "thailand", "thai-league-1", "2023-2024", 45283.54166666666, "Sukhothai", 2, "Khonkaen Utd.", 3, 1.76, "17/12/2023 11:42", 1.75, "23/12/2023 12:58", 3.79, "17/12/2023 11:42", 4.18, "23/12/2023 12:58", 4.03, "17/12/2023 11:42", 4.2, "23/12/2023 12:58", "https://www.betexplorer.com/football/thailand/thai-league-1/sukhothai-khonkaen-united/dl7ngNh1/"),
    @(111, "thailand", "thai-league-1", "2023-2024", 45283.58333333334, "Buriram", 4, "Uthai Thani", 0, 1.27, "17/12/2023 12:12", 1.43, "23/12/2023 13:52", 5.91, "17/12/2023 12:12", 5.32, "23/12/2023 13:57", 7.74, "17/12/2023 12:12", 6.17, "23/12/2023 13:57", "https://www.betexplorer.com/football/thailand/thai-league-1/buriram-united-f-c-uthai-thani/pAmsF4FK/"),
    @(112, "thailand", "thai-league-1", "2023-2024", 45284.47916666666, "Lamphun Warrior", 5, "Trat FC", 1, 1.78, "17/12/2023 13:12", 1.68, "24/12/2023 11:26", 3.71, "17/12/2023 13:12", 4.27, "24/12/2023 11:27", 4.03, "17/12/2023 13:12", 4.51, "24/12/2023 11:26", "https://www.betexplorer.com/football/thailand/thai-league-1/lamphun-warrior-trat-fc/0MDweqOl/"),
    @(113, "thailand", "thai-league-1", "2023-2024", 45284.5, "Chiangrai Utd", 1, "Prachuap", 1, 2.42, "17/12/2023 12:12", 2.39, "24/12/2023 11:48", 3.26, "17/12/2023 12:12", 3.28, "24/12/2023 11:45", 2.85, "17/12/2023 12:12", 3.12, "24/12/2023 11:52", "https://www.betexplorer.com/football/thailand/thai-league-1/chiangrai-utd-prachuap/2obxGp0E/"),
    @(114, "thailand", "thai-league-1", "2023-2024", 45284.52083333334, "Bangkok Utd", 1, "Police Tero", 1, 1.26, "17/12/2023 13:12", 1.17, "24/12/2023 12:04", 5.75, "17/12/2023 13:12", 7.66, "24/12/2023 12:04", 8.34, "17/12/2023 13:12", 15.41, "24/12/2023 12:04", "https://www.betexplorer.com/football/thailand/thai-league-1/bangkok-utd-police-tero/6aEZeP8r/"),
    @(115, "thailand", "thai-league-1", "2023-2024", 45284.54166666666, "Pathum United", 1, "Nakhon Pathom", 2, 1.32, "17/12/2023 13:12", 1.33, "24/12/2023 12:52", 5.19, "17/12/2023 13:12", 5.53, "24/12/2023 12:52", 7.21, "17/12/2023 13:12", 8.57, "24/12/2023 12:52", "https://www.betexplorer.com/football/thailand/thai-league-1/pathum-united-nakhon-pathom/vDCsf3we/"),
    @(116, "thailand", "thai-league-1", "2023-2024", 45285.54166666666, "Port MTI FC", 3, "Ratchaburi", 0, 1.54, "18/12/2023 21:12", 1.49, "25/12/2023 12:55", 4.41, "18/12/2023 21:12", 4.81, "25/12/2023 12:58", 4.93, "18/12/2023 21:12", 5.88, "25/12/2023 12:58", "https://www.betexplorer.com/football/thailand/thai-league-1/port-mti-fc-ratchaburi/lrfYGQp8/"),
    @(117, "thailand", "thai-league-1", "2023-2024", 45288.5, "Prachuap", 1, "Buriram", 2, 4.48, "25/12/2023 17:12", 4.56, "28/12/2023 11:59", 3.95, "25/12/2023 17:12", 3.94, "28/12/2023 11:59", 1.65, "25/12/2023 17:12", 1.74, "28/12/2023 11:59", "https://www.betexplorer.com/football/thailand/thai-league-1/prachuap-buriram-united-f-c/SbpqyUFI/"),
    @(118, "thailand", "thai-league-1", "2023-2024", 45288.54166666666, "Pathum United", 2, "Chiangrai Utd", 2, 1.52, "24/12/2023 13:12", 1.39, "28/12/2023 12:59", 4.2, "24/12/2023 13:12", 5.05, "28/12/2023 12:59", 5.49, "24/12/2023 13:12", 7.49, "28/12/2023 12:59", "https://www.betexplorer.com/football/thailand/thai-league-1/pathum-united-chiangrai-utd/hxxVbX7P/"),
    @(119, "thailand", "thai-league-1", "2023-2024", 45288.54166666666, "Chonburi", 0, "Bangkok Utd", 0, 3.66, "24/12/2023 12:42", 4.59, "28/12/2023 12:59", 3.65, "24/12/2023 12:42", 4.2, "28/12/2023 12:59", 1.88, "24/12/2023 12:42", 1.69, "28/12/2023 12:59", "https://www.betexplorer.com/football/thailand/thai-league-1/chonburi-bangkok-utd/bogUNJu2/")
)

$startRow = 110
for ($r = 0; $r -lt $newRows.Count; $r++) {
    $targetRow = $startRow + $r
    $data = $newRows[$r]

    # Copy the style/format (bold index column, bordered, date format, etc.)
    # from the previous last data row (109) down onto the new row, the same
    # way every other row in the sheet is formatted.
    $ws.Range("A109:V109").Copy()
    $ws.Range("A$targetRow" + ":V$targetRow").PasteSpecial(-4122)

    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value = $data[$c]
    }
}

$excel.CutCopyMode = 0

Write-Output "edit applied"
